$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 531, shifting rows 531-638 down to 532-639
$ws.Rows.Item(531).Insert()

# Populate the newly inserted row 531 with the new data
$ws.Cells.Item(531, 1).Value = 6
$ws.Cells.Item(531, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(531, 3).Value = "Metropolitana"

$d531 = Get-Date -Year 2023 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(531, 4).Value = $d531
$ws.Cells.Item(531, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(531, 5).Value = 13
$ws.Cells.Item(531, 6).Value = 100112032
$ws.Cells.Item(531, 7).Value = "Zapallo italiano"
$ws.Cells.Item(531, 8).Value = "Sin especificar"
$ws.Cells.Item(531, 9).Value = "Primera"
$ws.Cells.Item(531, 10).Value = 750
$ws.Cells.Item(531, 11).Value = 18000
$ws.Cells.Item(531, 12).Value = 20000
$ws.Cells.Item(531, 13).Value = 19200
$ws.Cells.Item(531, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(531, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(531, 16).Value = 384
$ws.Cells.Item(531, 17).Value = 50
$ws.Cells.Item(531, 18).Value = "Hortaliza"

# Update the sheet dimension to reflect the new last row
$ws.UsedRange | Out-Null
